# Apply the "产出平衡" updates: adjust the F column (粮/food) values for the
# "理想" (ideal) table (rows 18-23) on Sheet1. Dependent formulas (H, N, O
# columns) recalc automatically from these input changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 0.5
$ws.Range("F21").Value = 0.5
$ws.Range("F22").Value = 0.3
$ws.Range("F23").Value = 0.3

$excel.Calculate()

# Reflect the selection change seen in the saved file (user last clicked L29).
$ws.Range("L29").Select()
